$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.008.05"
$ws.Range("E2").Value = "  -0.36%  "
$ws.Range("D3").Value = "1.862.70"
$ws.Range("E3").Value = "  -0.89%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.30"
$ws.Range("E5").Value = "  -0.49%  "
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("E7").Value = "  +0.58%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3848"
$ws.Range("E8").Value = "  -0.51%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08257"
$ws.Range("E9").Value = "  -8.20%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.111"
$ws.Range("E10").Value = "  -1.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.54"
$ws.Range("E11").Value = "  -0.42%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.201"
$ws.Range("E12").Value = "  -2.60%  "
$ws.Range("E13").Value = "  -0.98%  "
$ws.Range("D14").Value = "1.864.17"
$ws.Range("E14").Value = "  -0.57%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.251"
$ws.Range("E15").Value = "  +0.19%  "
$ws.Range("E16").Value = "  +0.07%  "
$ws.Range("E17").Value = "  -1.03%  "
$ws.Range("E18").Value = "  -0.81%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06654"
$ws.Range("E19").Value = "  +0.51%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.70"
$ws.Range("E20").Value = "  -2.80%  "
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("E22").Value = "  -1.81%  "
$ws.Range("D23").Value = "28.032.81"
$ws.Range("E23").Value = "  -0.33%  "
$ws.Range("E24").Value = "  -3.20%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.231"
$ws.Range("E25").Value = "  -1.62%  "
$ws.Range("D26").Value = "2.076.26"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.512"
$ws.Range("E27").Value = "  -1.41%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "157.69"
$ws.Range("E28").Value = "  +0.53%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.49"
$ws.Range("E29").Value = "  -1.61%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "124.70"
$ws.Range("E30").Value = "  -1.92%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1063"
$ws.Range("E31").Value = "  +0.72%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.028"
$ws.Range("E32").Value = "  -3.28%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.940"
$ws.Range("E33").Value = "  +5.68%  "
$ws.Range("E34").Value = "  -0.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.377"
$ws.Range("E35").Value = "  -2.49%  "
$ws.Range("E36").Value = "  +0.18%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06497"
$ws.Range("E37").Value = "  -1.79%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2172"
$ws.Range("E38").Value = "  -0.62%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6583"
$ws.Range("E39").Value = "  +2.66%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.196"
$ws.Range("E40").Value = "  -1.31%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.993"
$ws.Range("E41").Value = "  +1.47%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.226"
$ws.Range("E42").Value = "  -4.79%  "
$ws.Range("E43").Value = "  -3.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6159"
$ws.Range("E44").Value = "  +1.92%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.01"
$ws.Range("E45").Value = "  -1.28%  "
$ws.Range("E46").Value = "  +0.51%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.668"
$ws.Range("E47").Value = "  -0.14%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.007"
$ws.Range("E48").Value = "  +0.37%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.217"
$ws.Range("E49").Value = "  -2.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "119.98"
$ws.Range("E50").Value = "  -1.16%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "78.86"
$ws.Range("E51").Value = "  -0.74%  "

Write-Output "Applied cryptos update"
